$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Colors (BGR ints as used by VBA/COM Interior.Color / Font.Color)
$Yellow = 65535      # FFFF00
$Red    = 255        # FF0000
$Green  = 5296274    # 92D050

# --- Column D: difficulty legend values (order matters for shared-string index) ---
$ws.Range("D2").Value = "Medium"
$ws.Range("D1").Value = "High"
$ws.Range("D3").Value = "Low"

# --- Row 1 (header row): red fill, centered for B1/D1, black automatic font ---
$ws.Range("A1").Interior.Color = $Red
$ws.Range("A1").Font.Name = "Calibri"
$ws.Range("A1").Font.Size = 11

$ws.Range("B1").Interior.Color = $Red
$ws.Range("B1").Font.Name = "Calibri"
$ws.Range("B1").Font.Size = 11
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108

$ws.Range("D1").Interior.Color = $Red
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4108

# --- Row 2: green ---
$ws.Range("A2").Interior.Color = $Green
$ws.Range("B2").Interior.Color = $Green
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("D2").Interior.Color = $Yellow
$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D2").VerticalAlignment = -4108

# --- Row 3: yellow ---
$ws.Range("A3").Interior.Color = $Yellow
$ws.Range("B3").Interior.Color = $Yellow
$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("D3").Interior.Color = $Green
$ws.Range("D3").HorizontalAlignment = -4108
$ws.Range("D3").VerticalAlignment = -4108

# --- Row 4: green ---
$ws.Range("A4").Interior.Color = $Green
$ws.Range("B4").Interior.Color = $Green
$ws.Range("B4").HorizontalAlignment = -4108
$ws.Range("B4").VerticalAlignment = -4108

# --- Row 5: red ---
$ws.Range("A5").Interior.Color = $Red
$ws.Range("B5").Interior.Color = $Red
$ws.Range("B5").HorizontalAlignment = -4108
$ws.Range("B5").VerticalAlignment = -4108

# --- Row 6: red ---
$ws.Range("A6").Interior.Color = $Red
$ws.Range("B6").Interior.Color = $Red
$ws.Range("B6").HorizontalAlignment = -4108
$ws.Range("B6").VerticalAlignment = -4108

# --- Row 7: yellow ---
$ws.Range("A7").Interior.Color = $Yellow
$ws.Range("B7").Interior.Color = $Yellow
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("B7").VerticalAlignment = -4108

# --- Row 8: yellow ---
$ws.Range("A8").Interior.Color = $Yellow
$ws.Range("B8").Interior.Color = $Yellow
$ws.Range("B8").HorizontalAlignment = -4108
$ws.Range("B8").VerticalAlignment = -4108

# --- Column B width / selection ---
$ws.Columns("B").ColumnWidth = 8.43
$ws.Range("A12").Select()
